$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Header row
$t.Cell(1, 2).Range.Text = "Categoria"
$t.Cell(1, 3).Range.Text = "Quantidade"

# Data rows (column 3 values updated with corrected results)
$t.Cell(2, 3).Range.Text = "365"
$t.Cell(3, 3).Range.Text = "224 (61.4)"
$t.Cell(4, 3).Range.Text = "141 (38.6)"
$t.Cell(5, 3).Range.Text = "63.56 (15.97)"
$t.Cell(6, 3).Range.Text = "1.63 (0.10)"
$t.Cell(7, 3).Range.Text = "75.51 (19.47)"
$t.Cell(8, 3).Range.Text = "28.40 (6.20)"
$t.Cell(9, 3).Range.Text = "140 (38.4)"
$t.Cell(10, 3).Range.Text = "225 (61.6)"
$t.Cell(11, 3).Range.Text = "357 (98.1)"
$t.Cell(12, 3).Range.Text = "7 ( 1.9)"
$t.Cell(13, 3).Range.Text = "325 (89.0)"
$t.Cell(14, 3).Range.Text = "40 (11.0)"
$t.Cell(15, 3).Range.Text = "322 (88.2)"
$t.Cell(16, 3).Range.Text = "43 (11.8)"
